$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite row 2 with the new student's first unit record
$ws.Range("A2").Value = "C026-01-0675/2020"
$ws.Range("B2").Value = "Clarence "
$ws.Range("C2").Value = "Gatama"
$ws.Range("D2").Value = "CCS4205"
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = "B"

# Overwrite row 3 with the same student's second unit record
$ws.Range("A3").Value = "C026-01-0675/2020"
$ws.Range("B3").Value = "Clarence "
$ws.Range("C3").Value = "Gatama"
$ws.Range("D3").Value = "CCS4204"
$ws.Range("E3").Value = 70
$ws.Range("F3").Value = "A"

# Remove the old fourth row entirely
$ws.Rows.Item(4).Delete()

# Re-balance column widths (shift C:F to take on B:E's widths)
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(5).ColumnWidth
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Match the resulting selection from the edit
$ws.Range("A4").Select()
